# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1.xml): rows 4,6,8,9,10
$wsExhibition.Range("F4").Value  = 175
$wsExhibition.Range("F6").Value  = 381
$wsExhibition.Range("F8").Value  = 2277
$wsExhibition.Range("F9").Value  = 385
$wsExhibition.Range("F10").Value = 5663

# 全部类型 sheet (sheet4.xml): rows 5,7,11,12,13
$wsAllTypes.Range("F5").Value  = 175
$wsAllTypes.Range("F7").Value  = 381
$wsAllTypes.Range("F11").Value = 2277
$wsAllTypes.Range("F12").Value = 385
$wsAllTypes.Range("F13").Value = 5663
